# Auto-generated-style edit script: update Leve market-price columns (H-N)
# across all 8 job sheets per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 958.0952
$ws.Range("I40").Value = 934.44446
$ws.Range("K40").Value = 934.44446
$ws.Range("M40").Value = -759.44446

$ws.Range("H76").Value = 3190.0732
$ws.Range("I76").Value = 3190.0732
$ws.Range("K76").Value = 3190.0732
$ws.Range("M76").Value = -2875.0732

$ws.Range("H79").Value = 3190.0732
$ws.Range("I79").Value = 3190.0732
$ws.Range("K79").Value = 3190.0732
$ws.Range("M79").Value = -2098.0732

$ws.Range("H112").Value = 76924580
$ws.Range("J112").Value = 76924580
$ws.Range("L112").Value = 230773740
$ws.Range("N112").Value = -230775956

$ws.Range("H127").Value = 1491.625
$ws.Range("J127").Value = 1966.6666
$ws.Range("L127").Value = 5899.9998
$ws.Range("N127").Value = -15819.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 400
$ws.Range("I16").Value = 400
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 400
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -113

$ws.Range("H32").Value = 13163.481
$ws.Range("I32").Value = 5830.244
$ws.Range("J32").Value = 36291.383
$ws.Range("K32").Value = 5830.244
$ws.Range("L32").Value = 36291.383
$ws.Range("M32").Value = -5543.244
$ws.Range("N32").Value = -36865.383

$ws.Range("H122").Value = 1735.2858
$ws.Range("I122").Value = 1374
$ws.Range("J122").Value = 2385.6
$ws.Range("K122").Value = 4122
$ws.Range("L122").Value = 7156.799999999999
$ws.Range("M122").Value = -1672
$ws.Range("N122").Value = -12056.8

$ws.Range("H132").Value = 2157857
$ws.Range("I132").Value = 2465928.5
$ws.Range("J132").Value = 1355.5
$ws.Range("K132").Value = 7397785.5
$ws.Range("L132").Value = 4066.5
$ws.Range("M132").Value = -7395255.5
$ws.Range("N132").Value = -9126.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 6682.8
$ws.Range("I75").Value = 3084.6667
$ws.Range("J75").Value = 12080
$ws.Range("K75").Value = 3084.6667
$ws.Range("L75").Value = 12080
$ws.Range("M75").Value = -2148.6667
$ws.Range("N75").Value = -13952

$ws.Range("H78").Value = 6682.8
$ws.Range("I78").Value = 3084.6667
$ws.Range("J78").Value = 12080
$ws.Range("K78").Value = 9254.000100000001
$ws.Range("L78").Value = 36240
$ws.Range("M78").Value = -4574.000100000001
$ws.Range("N78").Value = -45600

$ws.Range("H96").Value = 9483
$ws.Range("I96").Value = 8356.75
$ws.Range("J96").Value = 11735.5
$ws.Range("K96").Value = 8356.75
$ws.Range("L96").Value = 11735.5
$ws.Range("M96").Value = -5610.75
$ws.Range("N96").Value = -17227.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2009400
$ws.Range("I6").Value = 2511000
$ws.Range("K6").Value = 2511000
$ws.Range("M6").Value = -2510887

$ws.Range("H31").Value = 1292.1476
$ws.Range("I31").Value = 1029.8529
$ws.Range("J31").Value = 1622.4445
$ws.Range("K31").Value = 1029.8529
$ws.Range("L31").Value = 1622.4445
$ws.Range("M31").Value = -734.8529000000001
$ws.Range("N31").Value = -2212.4445

$ws.Range("H34").Value = 1292.1476
$ws.Range("I34").Value = 1029.8529
$ws.Range("J34").Value = 1622.4445
$ws.Range("K34").Value = 1029.8529
$ws.Range("L34").Value = 1622.4445
$ws.Range("M34").Value = -827.8529000000001
$ws.Range("N34").Value = -2026.4445

$ws.Range("H58").Value = 2842.2195
$ws.Range("I58").Value = 1457.85
$ws.Range("J58").Value = 4160.6665
$ws.Range("K58").Value = 1457.85
$ws.Range("L58").Value = 4160.6665
$ws.Range("M58").Value = -1254.85
$ws.Range("N58").Value = -4566.6665

$ws.Range("H94").Value = 1403.875
$ws.Range("I94").Value = 1153
$ws.Range("J94").Value = 1654.75
$ws.Range("K94").Value = 1153
$ws.Range("L94").Value = 1654.75
$ws.Range("M94").Value = -702
$ws.Range("N94").Value = -2556.75

$ws.Range("H132").Value = 3231.2307
$ws.Range("I132").Value = 2715
$ws.Range("K132").Value = 8145
$ws.Range("M132").Value = -5615

$ws.Range("H136").Value = 2842.2195
$ws.Range("I136").Value = 1457.85
$ws.Range("J136").Value = 4160.6665
$ws.Range("K136").Value = 4373.549999999999
$ws.Range("L136").Value = 12481.9995
$ws.Range("M136").Value = -1823.549999999999
$ws.Range("N136").Value = -17581.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 895.9878
$ws.Range("I68").Value = 677.2766
$ws.Range("J68").Value = 1189.6857
$ws.Range("K68").Value = 2031.8298
$ws.Range("L68").Value = 3569.0571
$ws.Range("M68").Value = -1220.8298
$ws.Range("N68").Value = -5191.0571

$ws.Range("H71").Value = 895.9878
$ws.Range("I71").Value = 677.2766
$ws.Range("J71").Value = 1189.6857
$ws.Range("K71").Value = 6095.4894
$ws.Range("L71").Value = 10707.1713
$ws.Range("M71").Value = -2039.4894
$ws.Range("N71").Value = -18819.1713

$ws.Range("H107").Value = 48385.047
$ws.Range("I107").Value = 30945.727
$ws.Range("J107").Value = 112329.22
$ws.Range("K107").Value = 92837.181
$ws.Range("L107").Value = 336987.66
$ws.Range("M107").Value = -90917.181
$ws.Range("N107").Value = -340827.66

$ws.Range("H131").Value = 1192143.6
$ws.Range("I131").Value = 979.3125
$ws.Range("J131").Value = 1472417.6
$ws.Range("K131").Value = 2937.9375
$ws.Range("L131").Value = 4417252.800000001
$ws.Range("M131").Value = 2102.0625
$ws.Range("N131").Value = -4427332.800000001

$ws.Range("H133").Value = 4821.905
$ws.Range("I133").Value = 1932.7273
$ws.Range("J133").Value = 8000
$ws.Range("K133").Value = 5798.1819
$ws.Range("L133").Value = 24000
$ws.Range("M133").Value = -738.1818999999996
$ws.Range("N133").Value = -34120

$ws.Range("H141").Value = 100003900
$ws.Range("I141").Value = 100003900
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 300011700
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -300006520

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 1550
$ws.Range("I14").Value = 1550
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1550
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -1382

$ws.Range("H102").Value = 1617.9412
$ws.Range("I102").Value = 1573.6
$ws.Range("J102").Value = 1681.2858
$ws.Range("K102").Value = 1573.6
$ws.Range("L102").Value = 1681.2858
$ws.Range("M102").Value = 48.40000000000009
$ws.Range("N102").Value = -4925.2858

$ws.Range("H113").Value = 1199.4286
$ws.Range("I113").Value = 1279.2
$ws.Range("K113").Value = 1279.2
$ws.Range("M113").Value = 890.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1667.6923
$ws.Range("I7").Value = 1737.3914
$ws.Range("J7").Value = 1133.3334
$ws.Range("K7").Value = 1737.3914
$ws.Range("L7").Value = 1133.3334
$ws.Range("M7").Value = -1625.3914
$ws.Range("N7").Value = -1357.3334

$ws.Range("H32").Value = 996
$ws.Range("I32").Value = 996
$ws.Range("K32").Value = 996
$ws.Range("M32").Value = -679

$ws.Range("H61").Value = 1961.8462
$ws.Range("I61").Value = 1309.4546
$ws.Range("K61").Value = 1309.4546
$ws.Range("M61").Value = -1107.4546

$ws.Range("H68").Value = 4129.4116
$ws.Range("I68").Value = 2050
$ws.Range("J68").Value = 4769.231
$ws.Range("K68").Value = 2050
$ws.Range("L68").Value = 4769.231
$ws.Range("M68").Value = -1301
$ws.Range("N68").Value = -6267.231

$ws.Range("H71").Value = 4129.4116
$ws.Range("I71").Value = 2050
$ws.Range("J71").Value = 4769.231
$ws.Range("K71").Value = 10250
$ws.Range("L71").Value = 23846.155
$ws.Range("M71").Value = -6506
$ws.Range("N71").Value = -31334.155

$ws.Range("H113").Value = 1961.8462
$ws.Range("I113").Value = 1309.4546
$ws.Range("K113").Value = 1309.4546
$ws.Range("M113").Value = 860.5454

$ws.Range("H122").Value = 16629.143
$ws.Range("I122").Value = 26126
$ws.Range("J122").Value = 3966.6667
$ws.Range("K122").Value = 78378
$ws.Range("L122").Value = 11900.0001
$ws.Range("M122").Value = -75928
$ws.Range("N122").Value = -16800.0001

$ws.Range("H126").Value = 1667.6923
$ws.Range("I126").Value = 1737.3914
$ws.Range("J126").Value = 1133.3334
$ws.Range("K126").Value = 5212.174199999999
$ws.Range("L126").Value = 3400.0002
$ws.Range("M126").Value = -2742.174199999999
$ws.Range("N126").Value = -8340.0002

$ws.Range("H132").Value = 5237.2334
$ws.Range("I132").Value = 5189.593
$ws.Range("J132").Value = 5666
$ws.Range("K132").Value = 15568.779
$ws.Range("L132").Value = 16998
$ws.Range("M132").Value = -13038.779
$ws.Range("N132").Value = -22058

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1073.6
$ws.Range("I100").Value = 1629.1428
$ws.Range("J100").Value = 587.5
$ws.Range("K100").Value = 3258.2856
$ws.Range("L100").Value = 1175
$ws.Range("M100").Value = -2717.2856
$ws.Range("N100").Value = -2257

$ws.Range("H107").Value = 318.45456
$ws.Range("I107").Value = 270.86667
$ws.Range("K107").Value = 812.60001
$ws.Range("M107").Value = 1107.39999

$ws.Range("H123").Value = 28612.154
$ws.Range("J123").Value = 28612.154
$ws.Range("L123").Value = 28612.154
$ws.Range("N123").Value = -38412.15399999999

$ws.Range("H126").Value = 1013.625
$ws.Range("I126").Value = 1120.8
$ws.Range("J126").Value = 835
$ws.Range("K126").Value = 3362.4
$ws.Range("L126").Value = 2505
$ws.Range("M126").Value = -892.3999999999996
$ws.Range("N126").Value = -7445

$ws.Range("H132").Value = 9040.565000000001
$ws.Range("I132").Value = 9040.565000000001
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 27121.695
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -24591.695

$ws.Range("H136").Value = 1153.7
$ws.Range("I136").Value = 1098.5294
$ws.Range("K136").Value = 3295.5882
$ws.Range("M136").Value = -745.5881999999997
